$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "23.993.81"
Set-TextValue "E2" "  -3.07%  "
Set-TextValue "D3" "1.631.10"
Set-TextValue "E3" "  -2.79%  "
Set-TextValue "D4" "1.006"
Set-TextValue "E4" "  +0.52%  "
Set-TextValue "B5" "USDC"
Set-TextValue "C5" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D5" "1.004"
Set-TextValue "E5" "  +0.21%  "
Set-TextValue "B6" "BNB"
Set-TextValue "C6" "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
Set-TextValue "D6" "307.36"
Set-TextValue "E6" "  -2.28%  "
Set-TextValue "D7" "0.3902"
Set-TextValue "E7" "  -0.38%  "
Set-TextValue "D8" "0.3843"
Set-TextValue "E8" "  -2.62%  "
Set-TextValue "D9" "1.004"
Set-TextValue "E9" "  +0.22%  "
Set-TextValue "D10" "49.99"
Set-TextValue "E10" "  -4.20%  "
Set-TextValue "D11" "1.357"
Set-TextValue "E11" "  -2.67%  "
Set-TextValue "D12" "0.08506"
Set-TextValue "E12" "  -1.57%  "
Set-TextValue "D13" "23.52"
Set-TextValue "E13" "  -6.90%  "
Set-TextValue "D14" "7.035"
Set-TextValue "E14" "  -3.98%  "
Set-TextValue "D15" "0.00001274"
Set-TextValue "E15" "  -3.29%  "
Set-TextValue "D16" "7.448"
Set-TextValue "E16" "  -4.04%  "
Set-TextValue "D17" "1.639.05"
Set-TextValue "E17" "  -1.34%  "
Set-TextValue "D18" "93.38"
Set-TextValue "E18" "  -0.48%  "
Set-TextValue "D19" "0.06911"
Set-TextValue "E19" "  -2.19%  "
Set-TextValue "D20" "20.13"
Set-TextValue "E20" "  -1.40%  "
Set-TextValue "D21" "6.876"
Set-TextValue "E21" "  -2.82%  "
Set-TextValue "D22" "1.003"
Set-TextValue "E22" "  +0.02%  "
Set-TextValue "D23" "13.51"
Set-TextValue "E23" "  -2.92%  "
Set-TextValue "D24" "24.030.25"
Set-TextValue "E24" "  -2.87%  "
Set-TextValue "D25" "2.414"
Set-TextValue "E25" "  +2.84%  "
Set-TextValue "D26" "2.832"
Set-TextValue "E26" "  +2.07%  "
Set-TextValue "D27" "22.10"
Set-TextValue "E27" "  -5.04%  "
Set-TextValue "D28" "157.66"
Set-TextValue "E28" "  -2.68%  "
Set-TextValue "D29" "139.55"
Set-TextValue "E29" "  -5.24%  "
Set-TextValue "D30" "5.246"
Set-TextValue "E30" "  -8.90%  "
Set-TextValue "D31" "7.706"
Set-TextValue "E31" "  -1.87%  "
Set-TextValue "D32" "2.444"
Set-TextValue "E32" "  -0.35%  "
Set-TextValue "D33" "1.818.20"
Set-TextValue "E33" "  -1.55%  "
Set-TextValue "D34" "0.08001"
Set-TextValue "E34" "  -4.98%  "
Set-TextValue "D35" "6.667"
Set-TextValue "E35" "  -3.78%  "
Set-TextValue "D36" "0.02878"
Set-TextValue "E36" "  -5.35%  "
Set-TextValue "D37" "0.9519"
Set-TextValue "E37" "  -4.33%  "
Set-TextValue "D38" "0.2669"
Set-TextValue "E38" "  -5.42%  "
Set-TextValue "D39" "0.09157"
Set-TextValue "E39" "  -3.66%  "
Set-TextValue "D40" "10.28"
Set-TextValue "E40" "  -2.90%  "
Set-TextValue "D41" "1.420"
Set-TextValue "E41" "  -8.53%  "
Set-TextValue "D42" "0.7447"
Set-TextValue "E42" "  -5.96%  "
Set-TextValue "D43" "12.98"
Set-TextValue "E43" "  -4.18%  "
Set-TextValue "D44" "16.06"
Set-TextValue "E44" "  -3.23%  "
Set-TextValue "D45" "0.6849"
Set-TextValue "E45" "  -4.07%  "
Set-TextValue "D46" "2.437"
Set-TextValue "E46" "  -4.92%  "
Set-TextValue "D47" "4.078"
Set-TextValue "E47" "  -2.83%  "
Set-TextValue "D48" "1.003"
Set-TextValue "E48" "  +0.22%  "
Set-TextValue "D49" "0.08294"
Set-TextValue "E49" "  -4.48%  "
Set-TextValue "B50" "Quant"
Set-TextValue "C50" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D50" "132.63"
Set-TextValue "E50" "  -3.84%  "
Set-TextValue "B51" "Flow"
Set-TextValue "C51" "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
Set-TextValue "D51" "1.250"
Set-TextValue "E51" "  -7.11%  "
